$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Metadata" ---
$meta = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date: 2021-12-16T17:36:56+00:00 -> 2022-01-21T20:46:54+00:00
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value (was blank) -> Alvearie Team
$meta.Range("B9").Value = "Alvearie Team"

# The old sheet had a duplicated "Contact" / "No display for ContactDetail" row
# (rows 10 and 11). Remove the duplicate row 11 entirely, then turn the
# remaining row 10 into the new "Jurisdiction" / "United States of America" row.
$meta.Rows.Item(11).Delete()
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# --- Sheet 2: "Elements" ---
$elements = $wb.Worksheets.Item("Elements")

# Root Extension row (row 2): Short / Definition text updated to be specific
# to this extension instead of the generic "Extension" / "An Extension".
$elements.Range("K2").Value = "Longterm Care Coverage Indicator"
$elements.Range("L2").Value = "Indicates whether the member or employee has long-term care benefit coverage"
